$wb = $excel.ActiveWorkbook

# Sheet ALC, row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 804.3036
$ws.Range("J17").Value = 812.43634
$ws.Range("L17").Value = 2437.30902
$ws.Range("N17").Value = -2773.30902

# Sheet ALC, row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1740
$ws.Range("I18").Value = 1740
$ws.Range("K18").Value = 1740
$ws.Range("M18").Value = -1456

# Sheet ALC, row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 773.8889
$ws.Range("I28").Value = 773.8889
$ws.Range("K28").Value = 773.8889
$ws.Range("M28").Value = -288.8889

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8250.166999999999
$ws.Range("I62").Value = 8647.235000000001
$ws.Range("J62").Value = 1500
$ws.Range("K62").Value = 8647.235000000001
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = -8023.235000000001
$ws.Range("N62").Value = -2748

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 8250.166999999999
$ws.Range("I65").Value = 8647.235000000001
$ws.Range("J65").Value = 1500
$ws.Range("K65").Value = 43236.175
$ws.Range("L65").Value = 7500
$ws.Range("M65").Value = -40116.175
$ws.Range("N65").Value = -13740

# Sheet ALC, row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2267.5
$ws.Range("I98").Value = 2335.8
$ws.Range("J98").Value = 1698.3334
$ws.Range("K98").Value = 2335.8
$ws.Range("L98").Value = 1698.3334
$ws.Range("M98").Value = -837.8000000000002
$ws.Range("N98").Value = -4694.3334

# Sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 7249
$ws.Range("I113").Value = 6999
$ws.Range("J113").Value = 7499
$ws.Range("K113").Value = 6999
$ws.Range("L113").Value = 7499
$ws.Range("M113").Value = -3745
$ws.Range("N113").Value = -14007

# Sheet ALC, row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2267.5
$ws.Range("I122").Value = 2335.8
$ws.Range("J122").Value = 1698.3334
$ws.Range("K122").Value = 7007.400000000001
$ws.Range("L122").Value = 5095.0002
$ws.Range("M122").Value = -4557.400000000001
$ws.Range("N122").Value = -9995.0002

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 21367.842
$ws.Range("I138").Value = 24797.297
$ws.Range("J138").Value = 12950.091
$ws.Range("K138").Value = 74391.891
$ws.Range("L138").Value = 38850.273
$ws.Range("M138").Value = -69251.891
$ws.Range("N138").Value = -49130.273

# Sheet ARM, row 31
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 3528.4
$ws.Range("I31").Value = 3528.4
$ws.Range("K31").Value = 3528.4
$ws.Range("M31").Value = -3234.4

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6259.1816
$ws.Range("I61").Value = 5807.1665
$ws.Range("J61").Value = 6801.6
$ws.Range("K61").Value = 5807.1665
$ws.Range("L61").Value = 6801.6
$ws.Range("M61").Value = -5595.1665
$ws.Range("N61").Value = -7225.6

# Sheet ARM, row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2196.7222
$ws.Range("I110").Value = 1253.75
$ws.Range("J110").Value = 4082.6667
$ws.Range("K110").Value = 1253.75
$ws.Range("L110").Value = 4082.6667
$ws.Range("M110").Value = 791.25
$ws.Range("N110").Value = -8172.6667

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 51303.81
$ws.Range("I132").Value = 65652
$ws.Range("K132").Value = 196956
$ws.Range("M132").Value = -194426

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6259.1816
$ws.Range("I136").Value = 5807.1665
$ws.Range("J136").Value = 6801.6
$ws.Range("K136").Value = 17421.4995
$ws.Range("L136").Value = 20404.8
$ws.Range("M136").Value = -14871.4995
$ws.Range("N136").Value = -25504.8

# Sheet BSM, row 12
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 500
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 500
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -836

# Sheet BSM, row 50
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 74999.5
$ws.Range("J50").Value = 74999.5
$ws.Range("L50").Value = 74999.5
$ws.Range("N50").Value = -76147.5

# Sheet CRP, row 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 13125844
$ws.Range("I6").Value = 13334259
$ws.Range("J6").Value = 12500600
$ws.Range("K6").Value = 13334259
$ws.Range("L6").Value = 12500600
$ws.Range("M6").Value = -13334146
$ws.Range("N6").Value = -12500826

# Sheet CRP, row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 68.375
$ws.Range("I7").Value = 46.2
$ws.Range("J7").Value = 105.333336
$ws.Range("K7").Value = 46.2
$ws.Range("L7").Value = 105.333336
$ws.Range("M7").Value = 66.8
$ws.Range("N7").Value = -331.333336

# Sheet CRP, row 17
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 715.2
$ws.Range("I17").Value = 715.2
$ws.Range("K17").Value = 715.2
$ws.Range("M17").Value = -541.2

# Sheet CRP, row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1015.9167
$ws.Range("I22").Value = 992.75
$ws.Range("K22").Value = 992.75
$ws.Range("M22").Value = -642.75

# Sheet CRP, row 25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 996.5
$ws.Range("I25").Value = 495.25
$ws.Range("J25").Value = 1999
$ws.Range("K25").Value = 495.25
$ws.Range("L25").Value = 1999
$ws.Range("M25").Value = -321.25
$ws.Range("N25").Value = -2347

# Sheet CRP, row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3184.6365
$ws.Range("I122").Value = 2754
$ws.Range("K122").Value = 8262
$ws.Range("M122").Value = -5812

# Sheet CUL, row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 414
$ws.Range("I26").Value = 460
$ws.Range("J26").Value = 368
$ws.Range("K26").Value = 1380
$ws.Range("L26").Value = 1104
$ws.Range("M26").Value = -1092
$ws.Range("N26").Value = -1680

# Sheet GSM, row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 157353.84
$ws.Range("I113").Value = 145745.14
$ws.Range("K113").Value = 145745.14
$ws.Range("M113").Value = -143575.14

# Sheet LTW, row 6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 68499
$ws.Range("J6").Value = 68499
$ws.Range("L6").Value = 68499
$ws.Range("N6").Value = -68723

# Sheet LTW, row 131
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Sheet WVR, row 51
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 29999
$ws.Range("I51").Value = 29999
$ws.Range("K51").Value = 29999
$ws.Range("M51").Value = -29489

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1371.2858
$ws.Range("I113").Value = 796
$ws.Range("J113").Value = 1601.4
$ws.Range("K113").Value = 2388
$ws.Range("L113").Value = 4804.200000000001
$ws.Range("M113").Value = -218
$ws.Range("N113").Value = -9144.200000000001

# Sheet WVR, row 127
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 99999
$ws.Range("J127").Value = 99999
$ws.Range("L127").Value = 99999
$ws.Range("N127").Value = -109919

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 44976.582
$ws.Range("J132").Value = 1196.5
$ws.Range("L132").Value = 3589.5
$ws.Range("N132").Value = -8649.5
